$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (matches source formatting)
$ws.Range('D2').Value = '66.946.00'
$ws.Range('E2').Value = '  +0.41%  '
$ws.Range('D3').Value = '3.497.76'
$ws.Range('E3').Value = '  +0.24%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '594.51'
$ws.Range('E5').Value = '  +0.43%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '172.37'
$ws.Range('E6').Value = '  +1.85%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  -0.36%  '
$ws.Range('E9').Value = '  +3.50%  '
$ws.Range('E10').Value = '  -1.27%  '
$ws.Range('E11').Value = '  -0.91%  '
$ws.Range('D12').Value = '4.100.83'
$ws.Range('E12').Value = '  +0.20%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '29.13'
$ws.Range('E14').Value = '  +3.46%  '
$ws.Range('D15').Value = '66.990.21'
$ws.Range('E15').Value = '  +0.43%  '
$ws.Range('E16').Value = '  +0.41%  '
$ws.Range('D17').Value = '3.492.86'
$ws.Range('E17').Value = '  +0.22%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.29'
$ws.Range('E18').Value = '  -0.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.29'
$ws.Range('E19').Value = '  +1.71%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '395.47'
$ws.Range('E20').Value = '  +0.70%  '
$ws.Range('E21').Value = '  +0.58%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '73.43'
$ws.Range('E22').Value = '  +0.43%  '
$ws.Range('E23').Value = '  +0.15%  '
$ws.Range('E24').Value = '  +0.36%  '
$ws.Range('E25').Value = '  -0.80%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.25'
$ws.Range('E26').Value = '  +0.20%  '
$ws.Range('E27').Value = '  -0.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.995'
$ws.Range('E28').Value = '  -0.62%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.19'
$ws.Range('E29').Value = '  -2.33%  '
$ws.Range('E30').Value = '  -2.50%  '
$ws.Range('E31').Value = '  -0.36%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '23.75'
$ws.Range('E32').Value = '  +0.91%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.37'
$ws.Range('E33').Value = '  -0.38%  '
$ws.Range('E34').Value = '  +0.55%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '162.90'
$ws.Range('E35').Value = '  +0.78%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.879'
$ws.Range('E36').Value = '  -2.50%  '
$ws.Range('E37').Value = '  -0.89%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.94'
$ws.Range('E38').Value = '  +2.79%  '
$ws.Range('E39').Value = '  +0.13%  '
$ws.Range('E40').Value = '  -0.52%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '27.20'
$ws.Range('E41').Value = '  +1.40%  '
$ws.Range('D42').Value = '2.836.41'
$ws.Range('E42').Value = '  +2.60%  '
$ws.Range('E43').Value = '  -1.20%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '42.82'
$ws.Range('E44').Value = '  -0.91%  '
$ws.Range('E45').Value = '  +2.41%  '
$ws.Range('E46').Value = '  -3.34%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '338.39'
$ws.Range('E47').Value = '  -2.31%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '34.75'
$ws.Range('E48').Value = '  +2.19%  '
$ws.Range('E49').Value = '  -1.41%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.843'
$ws.Range('E50').Value = '  -5.20%  '
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.105'
$ws.Range('E51').Value = '  -1.09%  '
